$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VehicleController")

$rows = @(45,48,49,50,51,52,53,59,61,63,68,70,71)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "Yes"
}

$ws.Range("E69").Select()
